$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 937
$ws.Range("I111").Value = 885.1111
$ws.Range("K111").Value = 2655.3333
$ws.Range("M111").Value = 411.6667000000002
$ws.Range("H134").Value = 39350
$ws.Range("J134").Value = 39350
$ws.Range("L134").Value = 39350
$ws.Range("N134").Value = -49490
$ws.Range("H137").Value = 2697.1462
$ws.Range("I137").Value = 2810.0908
$ws.Range("J137").Value = 2231.25
$ws.Range("K137").Value = 8430.2724
$ws.Range("L137").Value = 6693.75
$ws.Range("M137").Value = -5880.2724
$ws.Range("N137").Value = -11793.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2578.6924
$ws.Range("I2").Value = 2851
$ws.Range("K2").Value = 2851
$ws.Range("M2").Value = -2738
$ws.Range("H32").Value = 10304.56
$ws.Range("I32").Value = 6243.9414
$ws.Range("J32").Value = 36191
$ws.Range("K32").Value = 6243.9414
$ws.Range("L32").Value = 36191
$ws.Range("M32").Value = -5956.9414
$ws.Range("N32").Value = -36765
$ws.Range("H61").Value = 468361.97
$ws.Range("I61").Value = 528795.25
$ws.Range("K61").Value = 528795.25
$ws.Range("M61").Value = -528583.25
$ws.Range("H116").Value = 2578.6924
$ws.Range("I116").Value = 2851
$ws.Range("K116").Value = 2851
$ws.Range("M116").Value = -557
$ws.Range("H122").Value = 2880.8
$ws.Range("I122").Value = 1935.826
$ws.Range("J122").Value = 5985.7144
$ws.Range("K122").Value = 5807.478
$ws.Range("L122").Value = 17957.1432
$ws.Range("M122").Value = -3357.478
$ws.Range("N122").Value = -22857.1432
$ws.Range("H132").Value = 5122.72
$ws.Range("I132").Value = 8601.5
$ws.Range("J132").Value = 3485.647
$ws.Range("K132").Value = 25804.5
$ws.Range("L132").Value = 10456.941
$ws.Range("M132").Value = -23274.5
$ws.Range("N132").Value = -15516.941
$ws.Range("H136").Value = 468361.97
$ws.Range("I136").Value = 528795.25
$ws.Range("K136").Value = 1586385.75
$ws.Range("M136").Value = -1583835.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2578.6924
$ws.Range("I3").Value = 2851
$ws.Range("K3").Value = 2851
$ws.Range("M3").Value = -2737
$ws.Range("H107").Value = 2309.318
$ws.Range("I107").Value = 2090.111
$ws.Range("J107").Value = 3295.75
$ws.Range("K107").Value = 2090.111
$ws.Range("L107").Value = 3295.75
$ws.Range("M107").Value = -170.1109999999999
$ws.Range("N107").Value = -7135.75
$ws.Range("H134").Value = 2343.4614
$ws.Range("I134").Value = 1241.7142
$ws.Range("J134").Value = 4611.7646
$ws.Range("K134").Value = 3725.1426
$ws.Range("L134").Value = 13835.2938
$ws.Range("M134").Value = -1190.1426
$ws.Range("N134").Value = -18905.2938

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2616.884
$ws.Range("I31").Value = 1875.5778
$ws.Range("J31").Value = 4006.8333
$ws.Range("K31").Value = 1875.5778
$ws.Range("L31").Value = 4006.8333
$ws.Range("M31").Value = -1580.5778
$ws.Range("N31").Value = -4596.8333
$ws.Range("H34").Value = 2616.884
$ws.Range("I34").Value = 1875.5778
$ws.Range("J34").Value = 4006.8333
$ws.Range("K34").Value = 1875.5778
$ws.Range("L34").Value = 4006.8333
$ws.Range("M34").Value = -1673.5778
$ws.Range("N34").Value = -4410.8333
$ws.Range("H58").Value = 3863.4102
$ws.Range("I58").Value = 4316.393
$ws.Range("J58").Value = 2710.3635
$ws.Range("K58").Value = 4316.393
$ws.Range("L58").Value = 2710.3635
$ws.Range("M58").Value = -4113.393
$ws.Range("N58").Value = -3116.3635
$ws.Range("H132").Value = 2707.6287
$ws.Range("I132").Value = 1483.6428
$ws.Range("J132").Value = 3523.6191
$ws.Range("K132").Value = 4450.928400000001
$ws.Range("L132").Value = 10570.8573
$ws.Range("M132").Value = -1920.928400000001
$ws.Range("N132").Value = -15630.8573
$ws.Range("H136").Value = 3863.4102
$ws.Range("I136").Value = 4316.393
$ws.Range("J136").Value = 2710.3635
$ws.Range("K136").Value = 12949.179
$ws.Range("L136").Value = 8131.0905
$ws.Range("M136").Value = -10399.179
$ws.Range("N136").Value = -13231.0905

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2135.5
$ws.Range("I126").Value = 1880
$ws.Range("J126").Value = 2391
$ws.Range("K126").Value = 5640
$ws.Range("L126").Value = 7173
$ws.Range("M126").Value = -3170
$ws.Range("N126").Value = -12113
$ws.Range("H132").Value = 3899.1462
$ws.Range("I132").Value = 3944.5217
$ws.Range("K132").Value = 11833.5651
$ws.Range("M132").Value = -9303.5651

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2220.4211
$ws.Range("I7").Value = 1818.909
$ws.Range("J7").Value = 2772.5
$ws.Range("K7").Value = 1818.909
$ws.Range("L7").Value = 2772.5
$ws.Range("M7").Value = -1706.909
$ws.Range("N7").Value = -2996.5
$ws.Range("H126").Value = 2220.4211
$ws.Range("I126").Value = 1818.909
$ws.Range("J126").Value = 2772.5
$ws.Range("K126").Value = 5456.727000000001
$ws.Range("L126").Value = 8317.5
$ws.Range("M126").Value = -2986.727000000001
$ws.Range("N126").Value = -13257.5
$ws.Range("H132").Value = 11829.958
$ws.Range("I132").Value = 6357
$ws.Range("J132").Value = 14083.529
$ws.Range("K132").Value = 19071
$ws.Range("L132").Value = 42250.587
$ws.Range("M132").Value = -16541
$ws.Range("N132").Value = -47310.587
$ws.Range("H140").Value = 40402.715
$ws.Range("J140").Value = 40402.715
$ws.Range("L140").Value = 40402.715
$ws.Range("N140").Value = -50762.715
$ws.Range("H141").Value = 28443.334
$ws.Range("J141").Value = 28443.334
$ws.Range("L141").Value = 28443.334
$ws.Range("N141").Value = -38803.334

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1701.5
$ws.Range("I122").Value = 1665.1875
$ws.Range("J122").Value = 1798.3334
$ws.Range("K122").Value = 4995.5625
$ws.Range("L122").Value = 5395.0002
$ws.Range("M122").Value = -2545.5625
$ws.Range("N122").Value = -10295.0002
$ws.Range("H132").Value = 1682.3467
$ws.Range("I132").Value = 1006.11426
$ws.Range("J132").Value = 2274.05
$ws.Range("K132").Value = 3018.34278
$ws.Range("L132").Value = 6822.150000000001
$ws.Range("M132").Value = -488.3427799999999
$ws.Range("N132").Value = -11882.15
$ws.Range("H136").Value = 24444422
$ws.Range("I136").Value = 34518868
$ws.Range("J136").Value = 13207538
$ws.Range("K136").Value = 103556604
$ws.Range("L136").Value = 39622614
$ws.Range("M136").Value = -103554054
$ws.Range("N136").Value = -39627714
